$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update test case values (row 4 - boundary inputs)
$ws.Range("D4").Value = 188
$ws.Range("E4").Value = 150

# Row 6 (wheatInStore) - whole row bumped from 300 to 3000
$ws.Range("B6:F6").Value = 3000

# Row 10 (Return Value output) - acresOwned result updated
$ws.Range("E10").Value = 2950

# Update the selected cell to match the author's final selection
$ws.Range("E10").Select()
